# Update gh-pages to output generated at 456a3b4
#
# Both the "展览" sheet and the "全部类型" sheet previously led with two
# now-past 2024-08-03 events ("南宁·火影忍者only" and "南宁·蔚蓝档案only").
# The refreshed scrape drops those two rows (everything shifts up by two),
# renumbers the leading serial-number column, and bumps a few "想去人数"
# attendance counters that ticked up since the last scrape.

$wb = $excel.ActiveWorkbook

function Update-SheetForRefresh {
    param($SheetName)

    $ws = $wb.Worksheets.Item($SheetName)

    # Drop the two obsolete leading rows (2024-08-03 events); everything
    # below shifts up two rows automatically.
    $ws.Rows.Item(2).Delete()
    $ws.Rows.Item(2).Delete()

    # Renumber column A (serial number, 0-based: header row stays 0).
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}

Update-SheetForRefresh "展览"
Update-SheetForRefresh "全部类型"

# --- Refresh the "想去人数" (column F) counts that increased since the
# --- last scrape. Identify rows by their link (column H), which is
# --- unique per event and stable across the shift.
function Set-AttendanceByLink {
    param($SheetName, $Link, $NewCount)

    $ws = $wb.Worksheets.Item($SheetName)
    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    for ($r = 2; $r -le $lastRow; $r++) {
        if ($ws.Cells.Item($r, 8).Value2 -eq $Link) {
            $ws.Cells.Item($r, 6).Value = $NewCount
        }
    }
}

foreach ($sheetName in @("展览", "全部类型")) {
    # 南宁·熊喵M动漫嘉年华【免费】: 1059 -> 1066
    Set-AttendanceByLink $sheetName "https://show.bilibili.com/platform/detail.html?id=89145" 1066
    # 南宁·第二届北极光动漫展: 2458 -> 2467
    Set-AttendanceByLink $sheetName "https://show.bilibili.com/platform/detail.html?id=88276" 2467
    # 南宁·万圣漫控嘉年华10: 208 -> 209
    Set-AttendanceByLink $sheetName "https://show.bilibili.com/platform/detail.html?id=87820" 209
}
